$wb = $excel.ActiveWorkbook

# Remove the extra "Sheet 0" worksheet if it exists
foreach ($sheet in @($wb.Worksheets)) {
    if ($sheet.Name -eq "Sheet 0") {
        $sheet.Delete()
    }
}

$ws = $wb.Worksheets.Item("Sheet")

# Clear any existing content on the main sheet
$ws.Cells.Clear()

# Header row
$headers = @("Group 1", "Group 2", "Group 3")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Randomised data rows (words randomly distributed across the 3 groups)
$data = @(
    @("world", "are", "how"),
    @("you", "how", "hello"),
    @("you", "how", "you"),
    @("are", "are", "world"),
    @("world", "hello", "are"),
    @("hello", "world", "are"),
    @("today", "are", "how"),
    @("today", "you", "hello"),
    @("world", "how", "are"),
    @("hello", "you", "how"),
    @("how", "world", "today"),
    @("today", "are", "today"),
    @("today", "you", "how"),
    @("how", "world", "today"),
    @("are", "you", "you"),
    @("hello", "you", "today"),
    @("hello", "how", "world"),
    @("world", "are", "you"),
    @("world", "hello", "today"),
    @("today", "hello", "hello")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

$ws.Range("A1:E28").Select()
